# Insert a new weekly data row into the "Hortaliza, Feria Lagunitas de Puerto Montt - Ajo" sheet.
# The new row is inserted at row 481, pushing the existing rows 481:531 down to 482:532.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 481 (shifts row 481 and everything below it down by one).
$ws.Rows("481:481").Insert()

# Populate the newly inserted row 481 with the new weekly record.
$ws.Range("A481").Value = 4
$ws.Range("B481").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C481").Value = "Los Lagos"
$ws.Range("D481").Value = 45212
$ws.Range("E481").Value = 10
$ws.Range("F481").Value = 100112003
$ws.Range("G481").Value = "Ajo"
$ws.Range("H481").Value = "Chino"
$ws.Range("I481").Value = "Primera"
$ws.Range("J481").Value = 240
$ws.Range("K481").Value = 24000
$ws.Range("L481").Value = 25000
$ws.Range("M481").Value = 24500
$ws.Range("N481").Value = "$/caja 10 kilos"
$ws.Range("O481").Value = "China"
$ws.Range("P481").Value = 2450
$ws.Range("Q481").Value = 10
$ws.Range("R481").Value = "Hortaliza"
